$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.46594833333333
$ws.Range("H2").Value = 34.397845
$ws.Range("I2").Value = 0.0374233929424224
$ws.Range("J2").Value = 0.03742339294242241
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.36026066666667
$ws.Range("N2").Value = 49.080782
$ws.Range("O2").Value = 0.1040179164488296
$ws.Range("P2").Value = 0.1040179164488296
$ws.Range("Q2").Value = 187.5859035238655
$ws.Range("R2").Value = 1688.27313171479
$ws.Range("S2").Value = 0.003892703360316614
$ws.Range("T2").Value = 0.003892703360316615

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.46594833333333
$ws.Range("H3").Value = 34.397845
$ws.Range("I3").Value = 0.0374233929424224
$ws.Range("J3").Value = 0.03742339294242241
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.913432
$ws.Range("N3").Value = 83.740296
$ws.Range("O3").Value = 0.1774725413447623
$ws.Range("P3").Value = 0.1774725413447623
$ws.Range("Q3").Value = 320.0539691180133
$ws.Range("R3").Value = 2880.48572206212
$ws.Range("S3").Value = 0.006641624651235344
$ws.Range("T3").Value = 0.006641624651235346

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.46594833333333
$ws.Range("H4").Value = 34.397845
$ws.Range("I4").Value = 0.0374233929424224
$ws.Range("J4").Value = 0.03742339294242241
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.227296
$ws.Range("N4").Value = 42.68188799999999
$ws.Range("O4").Value = 0.09045660804390411
$ws.Range("P4").Value = 0.09045660804390411
$ws.Range("Q4").Value = 163.1294408590399
$ws.Range("R4").Value = 1468.16496773136
$ws.Range("S4").Value = 0.00338519318706571
$ws.Range("T4").Value = 0.003385193187065711

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.46594833333333
$ws.Range("H5").Value = 34.397845
$ws.Range("I5").Value = 0.0374233929424224
$ws.Range("J5").Value = 0.03742339294242241
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.972696
$ws.Range("N5").Value = 20.918088
$ws.Range("O5").Value = 0.04433213655506275
$ws.Range("P5").Value = 0.04433213655506275
$ws.Range("Q5").Value = 79.94857208003999
$ws.Range("R5").Value = 719.5371487203599
$ws.Range("S5").Value = 0.001659058966277241
$ws.Range("T5").Value = 0.001659058966277242

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.46594833333333
$ws.Range("H6").Value = 34.397845
$ws.Range("I6").Value = 0.0374233929424224
$ws.Range("J6").Value = 0.03742339294242241
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 48.33663566666667
$ws.Range("N6").Value = 145.009907
$ws.Range("O6").Value = 0.3073224951994155
$ws.Range("P6").Value = 0.3073224951994155
$ws.Range("Q6").Value = 554.2253671611571
$ws.Range("R6").Value = 4988.028304450415
$ws.Range("S6").Value = 0.01150105049789345
$ws.Range("T6").Value = 0.01150105049789345

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.46594833333333
$ws.Range("H7").Value = 34.397845
$ws.Range("I7").Value = 0.0374233929424224
$ws.Range("J7").Value = 0.03742339294242241
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 43.47278266666666
$ws.Range("N7").Value = 130.418348
$ws.Range("O7").Value = 0.2763983024080258
$ws.Range("P7").Value = 0.2763983024080257
$ws.Range("Q7").Value = 498.4566799622287
$ws.Range("R7").Value = 4486.110119660058
$ws.Range("S7").Value = 0.01034376227963404
$ws.Range("T7").Value = 0.01034376227963404

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 113.9117636666667
$ws.Range("H8").Value = 341.735291
$ws.Range("I8").Value = 0.3717934677996853
$ws.Range("J8").Value = 0.3717934677996854
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.36026066666667
$ws.Range("N8").Value = 49.080782
$ws.Range("O8").Value = 0.1040179164488296
$ws.Range("P8").Value = 0.1040179164488296
$ws.Range("Q8").Value = 1863.626146586395
$ws.Range("R8").Value = 16772.63531927756
$ws.Range("S8").Value = 0.0386731818698083
$ws.Range("T8").Value = 0.0386731818698083

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 113.9117636666667
$ws.Range("H9").Value = 341.735291
$ws.Range("I9").Value = 0.3717934677996853
$ws.Range("J9").Value = 0.3717934677996854
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.913432
$ws.Range("N9").Value = 83.740296
$ws.Range("O9").Value = 0.1774725413447623
$ws.Range("P9").Value = 0.1774725413447623
$ws.Range("Q9").Value = 3179.66826910957
$ws.Range("R9").Value = 28617.01442198613
$ws.Range("S9").Value = 0.06598313158579219
$ws.Range("T9").Value = 0.06598313158579219

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 113.9117636666667
$ws.Range("H10").Value = 341.735291
$ws.Range("I10").Value = 0.3717934677996853
$ws.Range("J10").Value = 0.3717934677996854
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.227296
$ws.Range("N10").Value = 42.68188799999999
$ws.Range("O10").Value = 0.09045660804390411
$ws.Range("P10").Value = 0.09045660804390411
$ws.Range("Q10").Value = 1620.656379567712
$ws.Range("R10").Value = 14585.9074161094
$ws.Range("S10").Value = 0.03363117599004001
$ws.Range("T10").Value = 0.03363117599004002

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 113.9117636666667
$ws.Range("H11").Value = 341.735291
$ws.Range("I11").Value = 0.3717934677996853
$ws.Range("J11").Value = 0.3717934677996854
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 6.972696
$ws.Range("N11").Value = 20.918088
$ws.Range("O11").Value = 0.04433213655506275
$ws.Range("P11").Value = 0.04433213655506275
$ws.Range("Q11").Value = 794.272098871512
$ws.Range("R11").Value = 7148.448889843608
$ws.Range("S11").Value = 0.01648239878477598
$ws.Range("T11").Value = 0.01648239878477598

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 113.9117636666667
$ws.Range("H12").Value = 341.735291
$ws.Range("I12").Value = 0.3717934677996853
$ws.Range("J12").Value = 0.3717934677996854
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 48.33663566666667
$ws.Range("N12").Value = 145.009907
$ws.Range("O12").Value = 0.3073224951994155
$ws.Range("P12").Value = 0.3073224951994155
$ws.Range("Q12").Value = 5506.111418503104
$ws.Range("R12").Value = 49555.00276652793
$ws.Range("S12").Value = 0.1142604962230428
$ws.Range("T12").Value = 0.1142604962230428

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 113.9117636666667
$ws.Range("H13").Value = 341.735291
$ws.Range("I13").Value = 0.3717934677996853
$ws.Range("J13").Value = 0.3717934677996854
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 43.47278266666666
$ws.Range("N13").Value = 130.418348
$ws.Range("O13").Value = 0.2763983024080258
$ws.Range("P13").Value = 0.2763983024080257
$ws.Range("Q13").Value = 4952.061345057696
$ws.Range("R13").Value = 44568.55210551926
$ws.Range("S13").Value = 0.102763083346226
$ws.Range("T13").Value = 0.102763083346226

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 61.749762
$ws.Range("H14").Value = 185.249286
$ws.Range("I14").Value = 0.2015433473312409
$ws.Range("J14").Value = 0.201543347331241
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 16.36026066666667
$ws.Range("N14").Value = 49.080782
$ws.Range("O14").Value = 0.1040179164488296
$ws.Range("P14").Value = 0.1040179164488296
$ws.Range("Q14").Value = 1010.242202424628
$ws.Range("R14").Value = 9092.179821821652
$ws.Range("S14").Value = 0.02096411906351847
$ws.Range("T14").Value = 0.02096411906351848

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 61.749762
$ws.Range("H15").Value = 185.249286
$ws.Range("I15").Value = 0.2015433473312409
$ws.Range("J15").Value = 0.201543347331241
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 27.913432
$ws.Range("N15").Value = 83.740296
$ws.Range("O15").Value = 0.1774725413447623
$ws.Range("P15").Value = 0.1774725413447623
$ws.Range("Q15").Value = 1723.647782603184
$ws.Range("R15").Value = 15512.83004342865
$ws.Range("S15").Value = 0.03576841004200544
$ws.Range("T15").Value = 0.03576841004200544

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 61.749762
$ws.Range("H16").Value = 185.249286
$ws.Range("I16").Value = 0.2015433473312409
$ws.Range("J16").Value = 0.201543347331241
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 14.227296
$ws.Range("N16").Value = 42.68188799999999
$ws.Range("O16").Value = 0.09045660804390411
$ws.Range("P16").Value = 0.09045660804390411
$ws.Range("Q16").Value = 878.5321419035517
$ws.Range("R16").Value = 7906.789277131966
$ws.Range("S16").Value = 0.01823092757339849
$ws.Range("T16").Value = 0.01823092757339849

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 61.749762
$ws.Range("H17").Value = 185.249286
$ws.Range("I17").Value = 0.2015433473312409
$ws.Range("J17").Value = 0.201543347331241
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.972696
$ws.Range("N17").Value = 20.918088
$ws.Range("O17").Value = 0.04433213655506275
$ws.Range("P17").Value = 0.04433213655506275
$ws.Range("Q17").Value = 430.562318498352
$ws.Range("R17").Value = 3875.060866485168
$ws.Range("S17").Value = 0.008934847195653016
$ws.Range("T17").Value = 0.008934847195653017

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 61.749762
$ws.Range("H18").Value = 185.249286
$ws.Range("I18").Value = 0.2015433473312409
$ws.Range("J18").Value = 0.201543347331241
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 48.33663566666667
$ws.Range("N18").Value = 145.009907
$ws.Range("O18").Value = 0.3073224951994155
$ws.Range("P18").Value = 0.3073224951994155
$ws.Range("Q18").Value = 2984.775748297378
$ws.Range("R18").Value = 26862.9817346764
$ws.Range("S18").Value = 0.06193880439267942
$ws.Range("T18").Value = 0.06193880439267942

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 61.749762
$ws.Range("H19").Value = 185.249286
$ws.Range("I19").Value = 0.2015433473312409
$ws.Range("J19").Value = 0.201543347331241
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 43.47278266666666
$ws.Range("N19").Value = 130.418348
$ws.Range("O19").Value = 0.2763983024080258
$ws.Range("P19").Value = 0.2763983024080257
$ws.Range("Q19").Value = 2684.433983144391
$ws.Range("R19").Value = 24159.90584829952
$ws.Range("S19").Value = 0.05570623906398611
$ws.Range("T19").Value = 0.0557062390639861

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 63.361323
$ws.Range("H20").Value = 190.083969
$ws.Range("I20").Value = 0.2068032768896493
$ws.Range("J20").Value = 0.2068032768896493
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 16.36026066666667
$ws.Range("N20").Value = 49.080782
$ws.Range("O20").Value = 0.1040179164488296
$ws.Range("P20").Value = 0.1040179164488296
$ws.Range("Q20").Value = 1036.607760464862
$ws.Range("R20").Value = 9329.469844183757
$ws.Range("S20").Value = 0.02151124597685172
$ws.Range("T20").Value = 0.02151124597685173

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 63.361323
$ws.Range("H21").Value = 190.083969
$ws.Range("I21").Value = 0.2068032768896493
$ws.Range("J21").Value = 0.2068032768896493
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 27.913432
$ws.Range("N21").Value = 83.740296
$ws.Range("O21").Value = 0.1774725413447623
$ws.Range("P21").Value = 0.1774725413447623
$ws.Range("Q21").Value = 1768.631980990536
$ws.Range("R21").Value = 15917.68782891482
$ws.Range("S21").Value = 0.03670190310803061
$ws.Range("T21").Value = 0.03670190310803061

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 63.361323
$ws.Range("H22").Value = 190.083969
$ws.Range("I22").Value = 0.2068032768896493
$ws.Range("J22").Value = 0.2068032768896493
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 14.227296
$ws.Range("N22").Value = 42.68188799999999
$ws.Range("O22").Value = 0.09045660804390411
$ws.Range("P22").Value = 0.09045660804390411
$ws.Range("Q22").Value = 901.4602972726078
$ws.Range("R22").Value = 8113.14267545347
$ws.Range("S22").Value = 0.01870672295980198
$ws.Range("T22").Value = 0.01870672295980198

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 63.361323
$ws.Range("H23").Value = 190.083969
$ws.Range("I23").Value = 0.2068032768896493
$ws.Range("J23").Value = 0.2068032768896493
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 6.972696
$ws.Range("N23").Value = 20.918088
$ws.Range("O23").Value = 0.04433213655506275
$ws.Range("P23").Value = 0.04433213655506275
$ws.Range("Q23").Value = 441.799243436808
$ws.Range("R23").Value = 3976.193190931272
$ws.Range("S23").Value = 0.009168031111106386
$ws.Range("T23").Value = 0.009168031111106388

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 63.361323
$ws.Range("H24").Value = 190.083969
$ws.Range("I24").Value = 0.2068032768896493
$ws.Range("J24").Value = 0.2068032768896493
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 48.33663566666667
$ws.Range("N24").Value = 145.009907
$ws.Range("O24").Value = 0.3073224951994155
$ws.Range("P24").Value = 0.3073224951994155
$ws.Range("Q24").Value = 3062.673185208987
$ws.Range("R24").Value = 27564.05866688088
$ws.Range("S24").Value = 0.06355529906914265
$ws.Range("T24").Value = 0.06355529906914265

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 63.361323
$ws.Range("H25").Value = 190.083969
$ws.Range("I25").Value = 0.2068032768896493
$ws.Range("J25").Value = 0.2068032768896493
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 43.47278266666666
$ws.Range("N25").Value = 130.418348
$ws.Range("O25").Value = 0.2763983024080258
$ws.Range("P25").Value = 0.2763983024080257
$ws.Range("Q25").Value = 2754.493024251467
$ws.Range("R25").Value = 24790.43721826321
$ws.Range("S25").Value = 0.05716007466471597
$ws.Range("T25").Value = 0.05716007466471597

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 40.18492533333333
$ws.Range("H26").Value = 120.554776
$ws.Range("I26").Value = 0.1311584709255395
$ws.Range("J26").Value = 0.1311584709255395
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 16.36026066666667
$ws.Range("N26").Value = 49.080782
$ws.Range("O26").Value = 0.1040179164488296
$ws.Range("P26").Value = 0.1040179164488296
$ws.Range("Q26").Value = 657.4358533238701
$ws.Range("R26").Value = 5916.922679914832
$ws.Range("S26").Value = 0.01364283087028902
$ws.Range("T26").Value = 0.01364283087028902

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 40.18492533333333
$ws.Range("H27").Value = 120.554776
$ws.Range("I27").Value = 0.1311584709255395
$ws.Range("J27").Value = 0.1311584709255395
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 27.913432
$ws.Range("N27").Value = 83.740296
$ws.Range("O27").Value = 0.1774725413447623
$ws.Range("P27").Value = 0.1774725413447623
$ws.Range("Q27").Value = 1121.699180717077
$ws.Range("R27").Value = 10095.2926264537
$ws.Range("S27").Value = 0.02327702715404861
$ws.Range("T27").Value = 0.02327702715404861

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 40.18492533333333
$ws.Range("H28").Value = 120.554776
$ws.Range("I28").Value = 0.1311584709255395
$ws.Range("J28").Value = 0.1311584709255395
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 14.227296
$ws.Range("N28").Value = 42.68188799999999
$ws.Range("O28").Value = 0.09045660804390411
$ws.Range("P28").Value = 0.09045660804390411
$ws.Range("Q28").Value = 571.7228274552318
$ws.Range("R28").Value = 5145.505447097087
$ws.Range("S28").Value = 0.01186415039614932
$ws.Range("T28").Value = 0.01186415039614932

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 40.18492533333333
$ws.Range("H29").Value = 120.554776
$ws.Range("I29").Value = 0.1311584709255395
$ws.Range("J29").Value = 0.1311584709255395
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 6.972696
$ws.Range("N29").Value = 20.918088
$ws.Range("O29").Value = 0.04433213655506275
$ws.Range("P29").Value = 0.04433213655506275
$ws.Range("Q29").Value = 280.197268132032
$ws.Range("R29").Value = 2521.775413188288
$ws.Range("S29").Value = 0.005814535243424244
$ws.Range("T29").Value = 0.005814535243424245

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 40.18492533333333
$ws.Range("H30").Value = 120.554776
$ws.Range("I30").Value = 0.1311584709255395
$ws.Range("J30").Value = 0.1311584709255395
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 48.33663566666667
$ws.Range("N30").Value = 145.009907
$ws.Range("O30").Value = 0.3073224951994155
$ws.Range("P30").Value = 0.3073224951994155
$ws.Range("Q30").Value = 1942.404095129537
$ws.Range("R30").Value = 17481.63685616583
$ws.Range("S30").Value = 0.04030794855137679
$ws.Range("T30").Value = 0.04030794855137679

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 40.18492533333333
$ws.Range("H31").Value = 120.554776
$ws.Range("I31").Value = 0.1311584709255395
$ws.Range("J31").Value = 0.1311584709255395
$ws.Range("K31").Value = 3
$ws.Range("M31").Value = 43.47278266666666
$ws.Range("N31").Value = 130.418348
$ws.Range("O31").Value = 0.2763983024080258
$ws.Range("P31").Value = 0.2763983024080257
$ws.Range("Q31").Value = 1746.950525492227
$ws.Range("R31").Value = 15722.55472943005
$ws.Range("S31").Value = 0.03625197871025151
$ws.Range("T31").Value = 0.03625197871025151

$ws.Range("E32").Value = 3
$ws.Range("G32").Value = 15.710799
$ws.Range("H32").Value = 47.132397
$ws.Range("I32").Value = 0.05127804411146253
$ws.Range("J32").Value = 0.05127804411146254
$ws.Range("K32").Value = 3
$ws.Range("M32").Value = 16.36026066666667
$ws.Range("N32").Value = 49.080782
$ws.Range("O32").Value = 0.1040179164488296
$ws.Range("P32").Value = 0.1040179164488296
$ws.Range("Q32").Value = 257.032766921606
$ws.Range("R32").Value = 2313.294902294454
$ws.Range("S32").Value = 0.00533383530804551
$ws.Range("T32").Value = 0.00533383530804551

$ws.Range("E33").Value = 3
$ws.Range("G33").Value = 15.710799
$ws.Range("H33").Value = 47.132397
$ws.Range("I33").Value = 0.05127804411146253
$ws.Range("J33").Value = 0.05127804411146254
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 27.913432
$ws.Range("N33").Value = 83.740296
$ws.Range("O33").Value = 0.1774725413447623
$ws.Range("P33").Value = 0.1774725413447623
$ws.Range("Q33").Value = 438.542319552168
$ws.Range("R33").Value = 3946.880875969512
$ws.Range("S33").Value = 0.009100444803650078
$ws.Range("T33").Value = 0.009100444803650078

$ws.Range("E34").Value = 3
$ws.Range("G34").Value = 15.710799
$ws.Range("H34").Value = 47.132397
$ws.Range("I34").Value = 0.05127804411146253
$ws.Range("J34").Value = 0.05127804411146254
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 14.227296
$ws.Range("N34").Value = 42.68188799999999
$ws.Range("O34").Value = 0.09045660804390411
$ws.Range("P34").Value = 0.09045660804390411
$ws.Range("Q34").Value = 223.5221877695039
$ws.Range("R34").Value = 2011.699689925535
$ws.Range("S34").Value = 0.004638437937448591
$ws.Range("T34").Value = 0.004638437937448591

$ws.Range("E35").Value = 3
$ws.Range("G35").Value = 15.710799
$ws.Range("H35").Value = 47.132397
$ws.Range("I35").Value = 0.05127804411146253
$ws.Range("J35").Value = 0.05127804411146254
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 6.972696
$ws.Range("N35").Value = 20.918088
$ws.Range("O35").Value = 0.04433213655506275
$ws.Range("P35").Value = 0.04433213655506275
$ws.Range("Q35").Value = 109.546625344104
$ws.Range("R35").Value = 985.919628096936
$ws.Range("S35").Value = 0.002273265253825888
$ws.Range("T35").Value = 0.002273265253825889

$ws.Range("E36").Value = 3
$ws.Range("G36").Value = 15.710799
$ws.Range("H36").Value = 47.132397
$ws.Range("I36").Value = 0.05127804411146253
$ws.Range("J36").Value = 0.05127804411146254
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 48.33663566666667
$ws.Range("N36").Value = 145.009907
$ws.Range("O36").Value = 0.3073224951994155
$ws.Range("P36").Value = 0.3073224951994155
$ws.Range("Q36").Value = 759.407167295231
$ws.Range("R36").Value = 6834.664505657079
$ws.Range("S36").Value = 0.01575889646528036
$ws.Range("T36").Value = 0.01575889646528036

$ws.Range("E37").Value = 3
$ws.Range("G37").Value = 15.710799
$ws.Range("H37").Value = 47.132397
$ws.Range("I37").Value = 0.05127804411146253
$ws.Range("J37").Value = 0.05127804411146254
$ws.Range("K37").Value = 3
$ws.Range("M37").Value = 43.47278266666666
$ws.Range("N37").Value = 130.418348
$ws.Range("O37").Value = 0.2763983024080258
$ws.Range("P37").Value = 0.2763983024080257
$ws.Range("Q37").Value = 682.9921504466839
$ws.Range("R37").Value = 6146.929354020155
$ws.Range("S37").Value = 0.0141731643432121
$ws.Range("T37").Value = 0.0141731643432121
